$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one
# day (45202 -> 45203) for every data row (rows 2 through 291).
for ($row = 2; $row -le 291; $row++) {
    $ws.Cells.Item($row, 3).Value = 45203
}
